$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 530
$ws.Range("I2").Value = 1500
$ws.Range("J2").Value = 6264
$ws.Range("K2").Value = 22
$ws.Range("L2").Value = 1670
$ws.Range("M2").Value = 89
$ws.Range("N2").Value = 1052
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = 86
$ws.Range("S2").Value = 654
$ws.Range("U2").Value = 86
$ws.Range("V2").Value = 9646
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 9635
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 128
$ws.Range("AA2").Value = 76
